$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 19, pushing the existing rows 19-30 down to 21-32.
$ws.Rows("19:20").Insert()

# --- New row 19: weekly "Primera" quote for 2022-04-13 ---
$ws.Cells.Item(19, 1).Value = 6
$ws.Cells.Item(19, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(19, 3).Value = "Metropolitana"
$ws.Cells.Item(19, 4).Value = 44664
$ws.Cells.Item(19, 5).Value = 13
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100101
$ws.Cells.Item(19, 8).Value = "Berries"
$ws.Cells.Item(19, 9).Value = 100101006
$ws.Cells.Item(19, 10).Value = "Higo"
$ws.Cells.Item(19, 11).Value = "Sin especificar"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 80
$ws.Cells.Item(19, 14).Value = 14000
$ws.Cells.Item(19, 15).Value = 14000
$ws.Cells.Item(19, 16).Value = 14000
$ws.Cells.Item(19, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(19, 18).Value = "Región Metropolitana"
$ws.Cells.Item(19, 19).Value = 2000
$ws.Cells.Item(19, 20).Value = 7

# --- New row 20: weekly "Segunda" quote for 2022-04-13 ---
$ws.Cells.Item(20, 1).Value = 6
$ws.Cells.Item(20, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(20, 3).Value = "Metropolitana"
$ws.Cells.Item(20, 4).Value = 44664
$ws.Cells.Item(20, 5).Value = 13
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100101
$ws.Cells.Item(20, 8).Value = "Berries"
$ws.Cells.Item(20, 9).Value = 100101006
$ws.Cells.Item(20, 10).Value = "Higo"
$ws.Cells.Item(20, 11).Value = "Sin especificar"
$ws.Cells.Item(20, 12).Value = "Segunda"
$ws.Cells.Item(20, 13).Value = 50
$ws.Cells.Item(20, 14).Value = 12000
$ws.Cells.Item(20, 15).Value = 12000
$ws.Cells.Item(20, 16).Value = 12000
$ws.Cells.Item(20, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(20, 18).Value = "Región Metropolitana"
$ws.Cells.Item(20, 19).Value = 1714
$ws.Cells.Item(20, 20).Value = 7
